$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell without Excel's automatic
# "this looks like a date" coercion (which would turn e.g. "1-2-2019" into a
# date serial number + a new date-formatted style). We stage the text in a
# scratch cell that is explicitly formatted as Text, copy it, and paste only
# the *value* into the destination cell - PasteSpecial(values) carries the
# "this is text" flag over without touching the destination cell's existing
# number format/style.
$scratch = $ws.Range("Z1")

function Set-TextValue($cellAddr, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)  # xlPasteValues
}

# New row 8 should look like the existing rows 5-7, so clone row 7's
# formatting down to row 8 before filling in the new values.
$ws.Range("A7:B7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

# Row 6: "2-2-2019" / "Error corrected in the file and changed " ->
#        "1-2-2019" / "Project error application properties"
Set-TextValue "A6" "1-2-2019"
$ws.Range("B6").Value = "Project error application properties"

# Row 7: "3-2-2019" / "Error in jdbc communication link failure" ->
#        "5-2-2019" / "Server port error xampp error resolved"
Set-TextValue "A7" "5-2-2019"
$ws.Range("B7").Value = "Server port error xampp error resolved"

# New row 8: "10-2-2019" / "Project completed successfully running on Postman "
Set-TextValue "A8" "10-2-2019"
$ws.Range("B8").Value = "Project completed successfully running on Postman "

# Clean up the scratch cell used for the text staging trick.
$scratch.Clear()

$ws.Range("B8").Select() | Out-Null
